# Atualizacao de bases das ligas, do dia: 19-04-2024 as 21:40
#
# This workbook holds "Australia ALeague" match/odds data, one row per game.
# The refreshed feed:
#   - removed the game that used to sit on row 153 (id 7127410), which
#     shifts every following row up by one;
#   - refreshed the odds figures for a handful of rows (both already
#     played games whose odds got recomputed, and upcoming games whose
#     odds moved), and in a couple of spots also swapped which of two
#     same-day fixtures occupies the earlier row.
#
# Strategy: delete the obsolete row first (so everything below shifts up
# and the shared-string / dimension bookkeeping is handled by Excel),
# then rewrite every cell of the affected rows with the refreshed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for game id 7127410 (25-04-2024 fixture that dropped out
# of the feed). Everything below (old rows 154-157) shifts up to 153-156.
$ws.Rows(153).Delete()

# Column letter -> index map used below: A=1 B=2 C=3 D=4 E=5 F=6 G=7 H=8
# I=9 J=10 K=11 L=12 M=13 N=14 O=15 P=16 Q=17 R=18 S=19 T=20 U=21 V=22
# W=23 X=24 Y=25 Z=26 AA=27 AB=28 AC=29

$rowsData = @(
  @{ row = 112; vals = @{
      1=110; 2=7127376; 3="Australia ALeague"; 4="Australia ALeague"; 5=45347.125;
      6="Newcastle Jets"; 7="Macarthur FC"; 8=2; 9=2; 10="D";
      11=1.95; 12=4; 13=3.4; 14=1.909; 15=4.2; 16=3.6; 17=-0.5;
      18=1.89; 19=2.01; 20=3.5; 21=1.95; 22=1.9;
      23=-1; 24=3.2; 25=-1; 26=-1; 27=1.01; 28=0.95; 29=-1
  } },
  @{ row = 113; vals = @{
      1=111; 2=7127379; 3="Australia ALeague"; 4="Australia ALeague"; 5=45347.125;
      6="Melbourne Victory"; 7="Central Coast Mariners"; 8=0; 9=1; 10="A";
      11=1.95; 12=3.6; 13=3.8; 14=1.909; 15=3.6; 16=4; 17=-0.5;
      18=1.9; 19=1.95; 20=2.75; 21=1.925; 22=1.925;
      23=-1; 24=-1; 25=3; 26=-1; 27=0.95; 28=-1; 29=0.925
  } },
  @{ row = 124; vals = @{
      1=122; 2=7128012; 3="Australia ALeague"; 4="Australia ALeague"; 5=45361.125;
      6="Macarthur FC"; 7="Central Coast Mariners"; 8=0; 9=3; 10="A";
      11=2.4; 12=3.5; 13=2.75; 14=3.4; 15=3.75; 16=2.05; 17=0.25;
      18=2.025; 19=1.825; 20=3; 21=2.05; 22=1.8;
      23=-1; 24=-1; 25=1.05; 26=-1; 27=0.825; 28=0; 29=0
  } },
  @{ row = 125; vals = @{
      1=123; 2=7127388; 3="Australia ALeague"; 4="Australia ALeague"; 5=45361.125;
      6="Sydney FC"; 7="Brisbane Roar"; 8=1; 9=1; 10="D";
      11=1.5; 12=5; 13=5; 14=1.533; 15=5.25; 16=5; 17=-1;
      18=1.8; 19=2.05; 20=3.5; 21=1.925; 22=1.925;
      23=-1; 24=4.25; 25=-1; 26=-1; 27=1.05; 28=-1; 29=0.925
  } },
  @{ row = 153; vals = @{
      1=151; 2=8096897; 3="Australia ALeague"; 4="Australia ALeague"; 5=45402.10416666666;
      6="Western Sydney Wanderers"; 7="Melbourne City";
      11=3.25; 12=3.8; 13=2; 14=3.3; 15=4; 16=2; 17=0.5;
      18=1.87; 19=2.03; 20=3.25; 21=1.925; 22=1.925;
      23=0; 24=0; 25=0; 26=0; 27=0
  } },
  @{ row = 154; vals = @{
      1=152; 2=7127411; 3="Australia ALeague"; 4="Australia ALeague"; 5=45402.1875;
      6="Melbourne Victory"; 7="Brisbane Roar";
      11=1.65; 12=4; 13=4.75; 14=1.6; 15=4.333; 16=5; 17=-1;
      18=2.05; 19=1.85; 20=3.25; 21=2.05; 22=1.8;
      23=0; 24=0; 25=0; 26=0; 27=0
  } },
  @{ row = 155; vals = @{
      1=153; 2=7127415; 3="Australia ALeague"; 4="Australia ALeague"; 5=45402.28125;
      6="Macarthur FC"; 7="Sydney FC";
      11=3.8; 12=4.2; 13=1.8; 14=4.333; 15=4.2; 16=1.666; 17=0.75;
      18=2.02; 19=1.88; 20=3.5; 21=1.925; 22=1.925;
      23=0; 24=0; 25=0; 26=0; 27=0
  } },
  @{ row = 156; vals = @{
      1=154; 2=7127414; 3="Australia ALeague"; 4="Australia ALeague"; 5=45403.16666666666;
      6="Perth Glory"; 7="Western United FC";
      11=2.4; 12=3.6; 13=2.625; 14=2.4; 15=3.75; 16=2.7; 17=0;
      18=1.84; 19=2.06; 20=3.5; 21=1.975; 22=1.875;
      23=0; 24=0; 25=0; 26=0; 27=0
  } }
)

foreach ($r in $rowsData) {
    foreach ($col in $r.vals.Keys) {
        $ws.Cells.Item($r.row, $col).Value = $r.vals[$col]
    }
}
